# Restore the "stack list" values in column A (tb3 working state).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 becomes the text label "stack list" using the sheet's default/base
# style (10pt), instead of the numeric "4" it previously held.
$ws.Range("A1").Value = "stack list"
$ws.Range("A1").Font.Size = 10

# A2:A5 get a small numeric stack list (1, 2, 3, 1) - previously blank.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 1

# Move/leave the active selection on A2, matching the saved view state.
$ws.Range("A2").Select()
